# Control Arm Naming Guide.xlsx - add a "QTY" column (H) to the Parts sheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Parts")

# New column header
$ws.Range("H1").Value = "QTY"

# Front/rear control-arm file-count summary row
$ws.Range("H4").Value = 24

# Front control arm part rows (10-16)
$ws.Range("H10").Value = 3
$ws.Range("H11").Value = 3
$ws.Range("H12").Value = 3
$ws.Range("H13").Value = 3
$ws.Range("H14").Value = 3
$ws.Range("H15").Value = 3
$ws.Range("H16").Value = 6

# Rear control arm part rows (23-30)
$ws.Range("H23").Value = 3
$ws.Range("H24").Value = 3
$ws.Range("H25").Value = 3
$ws.Range("H26").Value = 3
$ws.Range("H27").Value = 3
$ws.Range("H28").Value = 3
$ws.Range("H29").Value = 6
$ws.Range("H30").Value = 3

# Restore the view to scroll near the new data and select H10, matching
# the author's on-save cursor position.
$excel.ActiveWindow.ScrollRow = 16
$ws.Range("H10").Select()
